$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "SC 92" row (original row 28) and the "RM 232" row (original row 26).
# Deleting the higher-numbered row first keeps row 26's position valid.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# After the two row deletions, apply the individual "missing data" style cell edits.
$ws.Range("F2").Value2 = 18.03
$ws.Range("F3").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("F11").Value2 = 17.65
$ws.Range("F13").ClearContents()
$ws.Range("F21").Value2 = 16.58
$ws.Range("F25").ClearContents()
$ws.Range("B29").ClearContents()
$ws.Range("B33").Value2 = -19.5
$ws.Range("F33").Value2 = 17.53
